$d = $word.ActiveDocument

$pairs = @(
    @("90÷4=22, 2", "33÷2=16, 1"),
    @("79÷2=39, 1", "11÷6=1, 5"),
    @("19÷7=2, 5", "57÷9=6, 3"),
    @("75÷9=8, 3", "55÷8=6, 7"),
    @("98÷3=32, 2", "29÷7=4, 1"),
    @("39÷2=19, 1", "96÷9=10, 6"),
    @("41÷6=6, 5", "75÷3=25, 0"),
    @("41÷2=20, 1", "15÷6=2, 3"),
    @("68÷5=13, 3", "22÷2=11, 0"),
    @("79÷7=11, 2", "38÷8=4, 6"),
    @("82÷8=10, 2", "57÷4=14, 1"),
    @("23÷9=2, 5", "23÷7=3, 2"),
    @("74÷6=12, 2", "49÷6=8, 1"),
    @("96÷2=48, 0", "32÷9=3, 5"),
    @("34÷3=11, 1", "95÷3=31, 2"),
    @("13÷4=3, 1", "31÷8=3, 7"),
    @("56÷9=6, 2", "56÷2=28, 0"),
    @("13÷6=2, 1", "28÷9=3, 1"),
    @("19÷3=6, 1", "54÷8=6, 6"),
    @("87÷7=12, 3", "12÷5=2, 2"),
    @("80÷8=10, 0", "20÷7=2, 6"),
    @("91÷4=22, 3", "92÷8=11, 4"),
    @("43÷9=4, 7", "27÷8=3, 3"),
    @("91÷2=45, 1", "27÷7=3, 6"),
    @("16÷8=2, 0", "54÷7=7, 5")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
